$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(9966.2800000000007, 10051.719999999999, 18.84, 19,    $true,  0.85, 42613.766932870371, $false),
    @(9945.35,            9966.2800000000007, 18.93, 18.97, $true,  0.21, 42614.675358796296, $false),
    @(9854.85,            9945.35,            18.72, 18.89, $true,  0.91, 42615.752106481479, $false)
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 7).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($row, 8).Value = $vals[7]
}
